$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values that were re-pulled / recalculated
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 2
$ws.Range("F7").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("F12").Value = -4
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = -4
